## Update docs for version 1.1
## - Remove the "OutlierScope" row from the Config sheet (outlier detection
##   scope option was dropped).
## - Reword the IntraBatchMode explanation to describe the new Robust
##   (bisquare) linear regression method.
## - Reword the filter-bank explanation on the Clean & Explore sheet to
##   describe the new peak-wise summarisation options.
## - Leave the workbook with the "Clean & Explore" sheet active/selected.

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Config")
$wsClean  = $wb.Worksheets.Item("Clean & Explore")

# --- Config sheet -----------------------------------------------------
# Remove the whole "OutlierScope" row (row 4): Name/Options/Explanation
# all lived only in that row, so deleting it shifts every row below up
# by one and drops the now-unused option entirely.
$wsConfig.Activate()
$wsConfig.Rows.Item(4).Delete()

# The IntraBatchMode row is now row 7 (was row 8). Update its
# Explanation cell to describe the new Robust (bisquare) linear method.
$wsConfig.Range("C7").Value = 'Three correction modes. "Spline" is the default QCRSC algorithm that requires optimisation of the smoothing parameter. "Linear" is a simple Robust (bisquare) linear regression based on the QC values & requires no smoothing optimisation. "Mean" equalises the QC mean across batches & ignores within batch systematic change.'

# Update the on-screen scroll/selection to match the new layout.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$wsConfig.Range("F8").Select()

# --- Clean & Explore sheet --------------------------------------------
# Reword the filter-bank explanation (still row 3, column B).
$wsClean.Range("B3").Value = 'Should the filter bank calculate its peak-wise statistics across all batches or calulate each batch individually and then summarise based on the dropdown option (e.g. "Mean" = the mean of all the batch statistics & "Max" = the poorest statistic comparing all batches).'

# Make "Clean & Explore" the active/selected sheet & cell.
$wsClean.Activate()
$wsClean.Range("B3").Select()
